$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 (sheet 1) ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 07:18:13"
$ws1.Cells.Item(3,1).Value = "Total filas: 85"

$data1 = @(
    @(57, "07:18:13", "07:20", "10_OLMOS", 2, "LP1912"),
    @(58, "07:18:13", "07:21", "26_HERNANDEZ", 3, "LP1912"),
    @(59, "06:52:34", "07:23", "10_OLMOS", 31, "LP1912"),
    @(60, "06:52:34", "07:31", "11_ETCHEVERRY", 39, "LP1912"),
    @(61, "07:18:13", "07:32", "11_ETCHEVERRY", 14, "LP1912"),
    @(62, "07:18:13", "07:32", "84_COLONIA URQUIZA-ESC 49", 14, "LP1912"),
    @(63, "07:18:13", "07:35", "23_HERNANDEZ", 17, "LP1912"),
    @(64, "06:52:34", "07:36", "27_EL RETIRO", 44, "LP1912"),
    @(65, "07:18:13", "07:37", "27_EL RETIRO", 19, "LP1912"),
    @(66, "07:18:13", "07:39", "10_OLMOS", 21, "LP1912"),
    @(67, "05:54:50", "07:46", "16_SANTA ANA", 112, "LP1912"),
    @(68, "06:52:34", "07:47", "14_ABASTO", 55, "LP1912"),
    @(69, "07:18:13", "07:48", "14_ABASTO", 30, "LP1912"),
    @(70, "06:52:34", "07:51", "215D_EL PATO", 59, "LP1912"),
    @(71, "07:18:13", "07:52", "215D_EL PATO", 34, "LP1912"),
    @(72, "06:52:34", "07:58", "16_SANTA ANA", 66, "LP1912"),
    @(73, "07:18:13", "07:59", "16_SANTA ANA", 41, "LP1912"),
    @(74, "07:18:13", "08:00", "23_HERNANDEZ", 42, "LP1912"),
    @(75, "07:18:13", "08:05", "11_ETCHEVERRY", 47, "LP1912"),
    @(76, "06:24:49", "08:05", "23_HERNANDEZ", 101, "LP1912"),
    @(77, "06:52:34", "08:06", "23_HERNANDEZ", 74, "LP1912"),
    @(78, "07:18:13", "08:12", "15_ABASTO", 54, "LP1912"),
    @(79, "07:18:13", "08:21", "26_HERNANDEZ", 63, "LP1912"),
    @(80, "06:52:34", "08:22", "16_P MOR-SANTA ANA", 90, "LP1912"),
    @(81, "07:18:13", "08:23", "215B_EL PATO", 65, "LP1912"),
    @(82, "07:18:13", "08:23", "16_SANTA ANA", 65, "LP1912"),
    @(83, "07:18:13", "08:23", "16_P MOR-SANTA ANA", 65, "LP1912"),
    @(84, "07:18:13", "08:27", "84_COLONIA URQUIZA-ESC 49", 69, "LP1912"),
    @(85, "07:18:13", "08:42", "81_EL PELIGRO", 84, "LP1912"),
    @(86, "07:18:13", "08:44", "14_ABASTO", 86, "LP1912"),
    @(87, "07:18:13", "08:54", "17_ROMERO", 96, "LP1912"),
    @(88, "07:18:13", "09:02", "215A_EL PATO", 104, "LP1912"),
    @(89, "07:18:13", "09:11", "16_P MOR-SANTA ANA", 113, "LP1912"),
    @(90, "07:18:13", "09:17", "27_EL RETIRO", 119, "LP1912")
)
foreach ($row in $data1) {
    $r = $row[0]
    $ws1.Cells.Item($r,1).Value = $row[1]
    $ws1.Cells.Item($r,2).Value = $row[2]
    $ws1.Cells.Item($r,3).Value = $row[3]
    $ws1.Cells.Item($r,4).Value = $row[4]
    $ws1.Cells.Item($r,5).Value = $row[5]
}

# ---- Sheet: LP1912-215 (sheet 2) ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 07:18:13"
$ws2.Cells.Item(3,1).Value = "Total filas: 19"

$data2 = @(
    @(22, "07:18:13", "07:52", "215D_EL PATO", 34, "LP1912"),
    @(23, "07:18:13", "08:23", "215B_EL PATO", 65, "LP1912"),
    @(24, "07:18:13", "09:02", "215A_EL PATO", 104, "LP1912")
)
foreach ($row in $data2) {
    $r = $row[0]
    $ws2.Cells.Item($r,1).Value = $row[1]
    $ws2.Cells.Item($r,2).Value = $row[2]
    $ws2.Cells.Item($r,3).Value = $row[3]
    $ws2.Cells.Item($r,4).Value = $row[4]
    $ws2.Cells.Item($r,5).Value = $row[5]
}

# ---- Sheet: 6203-6173 (sheet 3) ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 07:18:13"
$ws3.Cells.Item(3,1).Value = "Total filas: 14"

$data3 = @(
    @(14, "07:18:13", "07:35", "215A_LA PLATA", 17, "L6173"),
    @(16, "07:18:13", "08:08", "215C_LA PLATA", 50, "L6203"),
    @(17, "06:52:34", "08:33", "215A_LA PLATA", 101, "L6173"),
    @(18, "07:18:13", "08:35", "215A_LA PLATA", 77, "L6173"),
    @(19, "07:18:13", "09:09", "215D_LA PLATA", 111, "L6203")
)
foreach ($row in $data3) {
    $r = $row[0]
    $ws3.Cells.Item($r,1).Value = $row[1]
    $ws3.Cells.Item($r,2).Value = $row[2]
    $ws3.Cells.Item($r,3).Value = $row[3]
    $ws3.Cells.Item($r,4).Value = $row[4]
    $ws3.Cells.Item($r,5).Value = $row[5]
}
